$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 21, shifting existing rows 21-40 down to 22-41
$ws.Rows("21:21").Insert()

# Populate the new row 21 with data (same constant fields as neighboring rows,
# with updated date / volume / price / origin values per the target diff)
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(22, 4).NumberFormat
$ws.Cells.Item(21, 4).Value = 44874
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112026
$ws.Cells.Item(21, 7).Value = "Haba"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 220
$ws.Cells.Item(21, 11).Value = 6000
$ws.Cells.Item(21, 12).Value = 6500
$ws.Cells.Item(21, 13).Value = 6227
$ws.Cells.Item(21, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value = "Región Metropolitana"
$ws.Cells.Item(21, 16).Value = 249
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"
